$wb = $excel.ActiveWorkbook

# Update values on sheet "展览" (rId1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 533
$ws1.Range("F3").Value = 33
$ws1.Range("F4").Value = 260
$ws1.Range("F5").Value = 3
$ws1.Range("F7").Value = 747

# Update values on sheet "全部类型" (rId4 / sheet4.xml) - mirrors the same data
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 533
$ws4.Range("F3").Value = 33
$ws4.Range("F4").Value = 260
$ws4.Range("F5").Value = 3
$ws4.Range("F7").Value = 747
